$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.351.05"
$ws.Range("E2").Value = "  +10.14%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.674.84"
$ws.Range("E3").Value = "  +5.67%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.44"
$ws.Range("E6").Value = "  +2.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3685"
$ws.Range("E7").Value = "  +2.33%  "

# Row 8
$ws.Range("E8").Value = "  +2.35%  "

# Row 9
$ws.Range("E9").Value = "  +14.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.157"
$ws.Range("E10").Value = "  +4.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07204"
$ws.Range("E11").Value = "  +4.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.127"
$ws.Range("E13").Value = "  +5.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.05"
$ws.Range("E14").Value = "  +4.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.710"
$ws.Range("E15").Value = "  +2.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.674.65"
$ws.Range("E16").Value = "  +5.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("E17").Value = "  +4.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  +0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.20"
$ws.Range("E20").Value = "  +6.00%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.46"
$ws.Range("E21").Value = "  +5.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.087"
$ws.Range("E22").Value = "  +3.08%  "

# Row 23
$ws.Range("E23").Value = "  +4.92%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.328.55"
$ws.Range("E24").Value = "  +9.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.434"
$ws.Range("E25").Value = "  +1.59%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.649"
$ws.Range("E26").Value = "  +7.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.06"
$ws.Range("E27").Value = "  +3.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("E28").Value = "  +1.64%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.859.29"
$ws.Range("E29").Value = "  +5.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.50"
$ws.Range("E30").Value = "  +5.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.275"
$ws.Range("E31").Value = "  +7.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.043"
$ws.Range("E32").Value = "  +1.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9674"
$ws.Range("E33").Value = "  +6.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08431"
$ws.Range("E34").Value = "  +3.87%  "

# Row 35
$ws.Range("E35").Value = "  +3.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.31"
$ws.Range("E36").Value = "  +6.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06390"
$ws.Range("E37").Value = "  +7.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.289"
$ws.Range("E38").Value = "  +4.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02313"
$ws.Range("E39").Value = "  +7.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.664"
$ws.Range("E40").Value = "  +4.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.231"
$ws.Range("E41").Value = "  +1.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2082"
$ws.Range("E42").Value = "  +5.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6064"
$ws.Range("E43").Value = "  +5.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.751"
$ws.Range("E45").Value = "  -0.61%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.90"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5859"
$ws.Range("E47").Value = "  +5.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.71"
$ws.Range("E48").Value = "  +0.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.009"
$ws.Range("E49").Value = "  +4.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07138"
$ws.Range("E50").Value = "  +6.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.71"
$ws.Range("E51").Value = "  +5.06%  "
